# Update "想去人数" (want-to-go count) values in column F across all four
# worksheets, per the upstream data refresh (gh-pages output regenerated
# at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 155
$ws.Range("F6").Value = 313
$ws.Range("F7").Value = 5547
$ws.Range("F8").Value = 65
$ws.Range("F9").Value = 52
$ws.Range("F10").Value = 3799
$ws.Range("F13").Value = 22
$ws.Range("F14").Value = 193
$ws.Range("F17").Value = 53
$ws.Range("F18").Value = 100
$ws.Range("F20").Value = 277
$ws.Range("F21").Value = 127
$ws.Range("F23").Value = 5215
$ws.Range("F25").Value = 2070
$ws.Range("F26").Value = 127
$ws.Range("F27").Value = 340
$ws.Range("F28").Value = 7731
$ws.Range("F31").Value = 2185
$ws.Range("F32").Value = 2156
$ws.Range("F34").Value = 156
$ws.Range("F35").Value = 1180
$ws.Range("F37").Value = 17
$ws.Range("F38").Value = 259
$ws.Range("F40").Value = 243
$ws.Range("F41").Value = 12
$ws.Range("F45").Value = 1321
$ws.Range("F46").Value = 2035
$ws.Range("F47").Value = 118
$ws.Range("F48").Value = 210
$ws.Range("F49").Value = 1211

# --- Sheet "演出" (sheet2) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 15
$ws.Range("F20").Value = 1

# --- Sheet "本地生活" (sheet3) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 551
$ws.Range("F3").Value = 726

# --- Sheet "全部类型" (sheet4) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 155
$ws.Range("F6").Value = 551
$ws.Range("F7").Value = 726
$ws.Range("F8").Value = 313
$ws.Range("F9").Value = 5547
$ws.Range("F10").Value = 3799
$ws.Range("F13").Value = 22
$ws.Range("F14").Value = 193
$ws.Range("F16").Value = 53
$ws.Range("F17").Value = 100
$ws.Range("F18").Value = 15
$ws.Range("F20").Value = 277
$ws.Range("F22").Value = 127
$ws.Range("F24").Value = 5215
$ws.Range("F26").Value = 2070
$ws.Range("F27").Value = 127
$ws.Range("F28").Value = 340
$ws.Range("F29").Value = 7731
$ws.Range("F32").Value = 2185
$ws.Range("F33").Value = 2156
$ws.Range("F35").Value = 156
$ws.Range("F36").Value = 1180
$ws.Range("F37").Value = 259
$ws.Range("F38").Value = 243
$ws.Range("F39").Value = 12
$ws.Range("F43").Value = 1321
$ws.Range("F45").Value = 2035
$ws.Range("F46").Value = 118
$ws.Range("F48").Value = 210
$ws.Range("F49").Value = 1211
